# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet. It gets
#    the same per-fund holdings table layout as the other quarterly sheets,
#    so clone "2021-Q4" (which already carries the right headers/formatting)
#    and just overwrite its data row.
# 2. Update the "总计" (totals) sheet: insert a new first data row for
#    "2022-Q1" and keep the running index column (column A) sequential.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "2022-Q1" worksheet, positioned right before "总计" ---
$wb.Worksheets.Item("2021-Q4").Copy($wb.Worksheets.Item("总计"))
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "320017"
$newSheet.Range("C2").Value = "诺安全球收益不动产(QDII)"
$newSheet.Range("D2").Value = "0.29"
$newSheet.Range("E2").Value = "93.32"
$newSheet.Range("F2").Value = "5.76"
$newSheet.Range("G2").Value = "0.0167"
$newSheet.Range("B2:G2").Style = "Normal"

$newSheet.Range("H2").Value = 6

# --- 2. Update the "总计" sheet: insert the 2022-Q1 row at the top ---
# (Re-fetch by name: the sheet collection shifted after the Copy() above, so
#  any reference captured before that point would now point at the wrong
#  sheet.)
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows(2).Insert()

# Give the new A2 cell the same formatting as the index column below it.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A2").Value = 0

$totalSheet.Range("B2").NumberFormat = "@"
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("B2").Style = "Normal"

$totalSheet.Range("C2").Value = 1
$totalSheet.Range("C2").Style = "Normal"

$totalSheet.Range("D2").Value = 0.02
$totalSheet.Range("D2").Style = "Normal"

# Keep the running index in column A sequential (0,1,2,3,4) for the rows
# that were pushed down by the insert.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

# Restore the originally active sheet/tab (the edit should not change which
# sheet is selected).
$wb.Worksheets.Item("2020-Q4").Activate()
